# Apply "Create default order for menu" change:
# - Insert a new "Gender" cut row into the Lookups config table (row 3),
#   pushing cuts_historical down by one row.
# - Re-order the "cuts" default-menu column (F) so it no longer mirrors
#   the config table order exactly (Region, Grade, Gender) and append a
#   new historical marker value "R".
# - Add two new (currently empty) defined names cuts_2 / cuts_3 reserved
#   for future cuts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# --- Update the Lookups sheet cell values ---

# Swap F1/F2 default-menu order relative to the config table.
$ws.Range("F1").Value = "Region"
$ws.Range("F2").Value = "Grade"

# New row 3: a "Gender" cut definition (config columns A-E) plus its
# default-menu entry in column F.
$ws.Range("A3").Value = "Gender"
$ws.Range("B3").Value = "static"
$ws.Range("C3").Value = "Gender"
$ws.Range("D3").Value = "Region"
$ws.Range("E3").Value = "Corps"
$ws.Range("F3").Value = "Gender"

# Row 4: the historical cut marker moves down to F4 with a new value.
$ws.Range("F4").Value = "R"

# --- Update the defined names ---

$wb.Names.Item("cuts_config").RefersTo = "=Lookups!`$A`$1:`$E`$3"
$wb.Names.Item("cuts").RefersTo = "=Lookups!`$F`$1:`$F`$3"
$wb.Names.Item("cuts_historical").RefersTo = "=Lookups!`$F`$4:`$F`$4"

$wb.Names.Add("cuts_2", "=Lookups!`$F`$5:`$F`$4")
$wb.Names.Add("cuts_3", "=Lookups!`$F`$5:`$F`$4")
